$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 95, pushing existing rows 95-147 down to 97-149.
$ws.Rows("95:96").Insert()

# --- New row 95 ---
$ws.Cells.Item(95, 1).Value = 9
$ws.Cells.Item(95, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(95, 3).Value = "Metropolitana"
$ws.Cells.Item(95, 4).Value = 44529
$ws.Cells.Item(95, 5).Value = 13
$ws.Cells.Item(95, 6).Value = 100112003
$ws.Cells.Item(95, 7).Value = "Ajo"
$ws.Cells.Item(95, 8).Value = "Chino"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 520
$ws.Cells.Item(95, 11).Value = 17500
$ws.Cells.Item(95, 12).Value = 18000
$ws.Cells.Item(95, 13).Value = 17750
$ws.Cells.Item(95, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(95, 15).Value = "China"
$ws.Cells.Item(95, 16).Value = 1775
$ws.Cells.Item(95, 17).Value = 10
$ws.Cells.Item(95, 18).Value = "Hortaliza"

# --- New row 96 ---
$ws.Cells.Item(96, 1).Value = 9
$ws.Cells.Item(96, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(96, 3).Value = "Metropolitana"
$ws.Cells.Item(96, 4).Value = 44529
$ws.Cells.Item(96, 5).Value = 13
$ws.Cells.Item(96, 6).Value = 100112003
$ws.Cells.Item(96, 7).Value = "Ajo"
$ws.Cells.Item(96, 8).Value = "Chino"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 610
$ws.Cells.Item(96, 11).Value = 17000
$ws.Cells.Item(96, 12).Value = 17500
$ws.Cells.Item(96, 13).Value = 17250
$ws.Cells.Item(96, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(96, 15).Value = "China"
$ws.Cells.Item(96, 16).Value = 1725
$ws.Cells.Item(96, 17).Value = 10
$ws.Cells.Item(96, 18).Value = "Hortaliza"
